$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing "sum" header (G1) to the new H1 header cell
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Set the new header text and value
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
